# Apply "Penalty Reward System" forecast update:
#  - Shift each week's Week_Start_Date forward by one week
#  - Update MyForecast (column D) values on "Forecast Comparison" sheet
#  - Update derived metrics on the "Summary" sheet

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# Make sure the date-like text we write stays as literal text (not auto
# converted to an Excel date serial number) by pre-formatting the ranges
# as Text before assigning values.
$ws1.Range("B2:B17").NumberFormat = "@"
$ws2.Range("B2:B15").NumberFormat = "@"

# --- "Forecast Comparison" sheet: Week_Start_Date (col B) ---
$ws1.Range("B2").Value  = "2025-01-12"
$ws1.Range("B3").Value  = "2025-01-19"
$ws1.Range("B4").Value  = "2025-01-26"
$ws1.Range("B5").Value  = "2025-02-02"
$ws1.Range("B6").Value  = "2025-02-09"
$ws1.Range("B7").Value  = "2025-02-16"
$ws1.Range("B8").Value  = "2025-02-23"
$ws1.Range("B9").Value  = "2025-03-02"
$ws1.Range("B10").Value = "2025-03-09"
$ws1.Range("B11").Value = "2025-03-16"
$ws1.Range("B12").Value = "2025-03-23"
$ws1.Range("B13").Value = "2025-03-30"
$ws1.Range("B14").Value = "2025-04-06"
$ws1.Range("B15").Value = "2025-04-13"
$ws1.Range("B16").Value = "2025-04-20"
$ws1.Range("B17").Value = "2025-04-27"

# --- "Forecast Comparison" sheet: MyForecast (col D) ---
$ws1.Range("D2").Value  = 2
$ws1.Range("D3").Value  = 1
$ws1.Range("D4").Value  = 2
$ws1.Range("D5").Value  = 2
$ws1.Range("D6").Value  = 2
$ws1.Range("D7").Value  = 2
$ws1.Range("D8").Value  = 2
$ws1.Range("D9").Value  = 2
$ws1.Range("D10").Value = 2
$ws1.Range("D11").Value = 2
$ws1.Range("D12").Value = 2
$ws1.Range("D13").Value = 2
$ws1.Range("D14").Value = 2
$ws1.Range("D15").Value = 2
$ws1.Range("D16").Value = 2
$ws1.Range("D17").Value = 2

# --- "Summary" sheet updates ---
$ws2.Range("B2").Value  = "2023-01-01 to 2025-01-05"
$ws2.Range("B4").Value  = "5"
$ws2.Range("B7").Value  = "2"
$ws2.Range("B8").Value  = "47 units"
$ws2.Range("B9").Value  = "30"
$ws2.Range("B10").Value = "13"
$ws2.Range("B11").Value = "6"
$ws2.Range("B12").Value = "2"
$ws2.Range("B13").Value = "2025-03-30"
$ws2.Range("B14").Value = "1"
$ws2.Range("B15").Value = "2025-01-19"
